# anzahl Schüler und Gebäudegrösse (Schulhausgrösse) - add two new data rows
# at the bottom of the sheet (rows 20 & 21), matching the style of the
# existing "label" rows above them, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: "anzahl Schüler" -------------------------------------------
# A2 carries the bold/Verdana "section label" style (cellXf index 2) that
# the other main labels in column A use - copy it onto A20 before writing
# the new header text.
$ws.Range("A2").Copy()
$ws.Range("A20").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A20").Value = "anzahl Schüler"

$ws.Range("B20").Value = 80
$ws.Range("C20").Value = 120
$ws.Range("D20").Value = 200
$ws.Range("E20").Value = 180
$ws.Range("F20").Value = 200

# --- Row 21: "Gebäudegrösse in Quadratmeter" ----------------------------
# A4 carries the plain wrap-text/shaded style (cellXf index 1) used by the
# spacer rows - copy it onto A21 before writing the new header text.
$ws.Range("A4").Copy()
$ws.Range("A21").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A21").Value = "Gebäudegrösse in Quadratmeter"

$ws.Range("B21").Value = 170
$ws.Range("C21").Value = 200
$ws.Range("D21").Value = 250
$ws.Range("E21").Value = 240
$ws.Range("F21").Value = 300

# Row-height tweaks that accompanied the edit in the source workbook.
$ws.Rows.Item(3).RowHeight = 22.5
$ws.Rows.Item(10).RowHeight = 32.45

# Leave the clipboard clean and land the selection where the author's
# workbook ended up after entering the new data.
$excel.CutCopyMode = 0
$ws.Range("B24").Select() | Out-Null
